$d = $word.ActiveDocument

# Locate the "Edison Achalma" paragraph styled as "Author" (the author byline
# right under the title), so we can add the affiliation paragraph after it.
$rng = $d.Content
$found = $rng.Find.Execute("Edison Achalma", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $authorPara = $rng.Paragraphs(1)

    # Insert a new paragraph right after the author's paragraph mark.
    $insertPoint = $d.Range($authorPara.Range.End, $authorPara.Range.End)
    $insertPoint.InsertParagraphAfter()

    $newPara = $authorPara.Next()
    $newPara.Style = "Author"
    $newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
}
